$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Teh kunyit", 24),
    @("Teh Kombucha", 18),
    @("Teh barley", 19),
    @("Teh dandelion", 10),
    @("Teh yerba mate", 34),
    @("Teh krisan", 43),
    @("Teh moringa", 23)
)

$row = 15
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B22").Select()
